# Revert "removing _new.xlsx in the databases"
#
# The HEATING and COOLING sheets had gained three extra columns
# (primary_components, secondary_components, tertiary_components,
# shared-string values "-", "BO2", "BO4", ... etc.) - this change removes
# those columns again, restoring the original (pre-addition) layout that
# HOT_WATER and ELECTRICITY already have.

$wb = $excel.ActiveWorkbook

# --- HEATING: drop columns C:E (primary/secondary/tertiary_components) ---
$wsHeating = $wb.Worksheets.Item("HEATING")
$wsHeating.Range("C1:E1").EntireColumn.Delete()

# --- COOLING: drop columns C:E (primary/secondary/tertiary_components) ---
$wsCooling = $wb.Worksheets.Item("COOLING")
$wsCooling.Range("C1:E1").EntireColumn.Delete()

# --- restore the view/selection state recorded in the workbook ---

# HOT_WATER: selection/layout untouched by the column edit.
$wsHotWater = $wb.Worksheets.Item("HOT_WATER")
$wsHotWater.Range("C1:C1048576").Select()

# COOLING: selection moves off the (now gone) component columns.
$wsCooling.Range("D10").Select()

# ELECTRICITY: keep its existing selection, just lose the "active tab" flag
# (HEATING becomes the active tab below).
$wsElectricity = $wb.Worksheets.Item("ELECTRICITY")
$wsElectricity.Range("B7").Select()

# HEATING becomes the active sheet/tab, with a new selected cell.
$wsHeating.Activate()
$wsHeating.Range("E12").Select()
